$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.657.18"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "1.902.52"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D4").Value = "'1.026"
$ws.Range("E4").Value = "  +2.16%  "
$ws.Range("D5").Value = "'320.06"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").Value = "'1.028"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("D7").Value = "'0.5185"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("D8").Value = "'0.3956"
$ws.Range("E8").Value = "  +3.37%  "
$ws.Range("D9").Value = "'0.08362"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "'1.135"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").Value = "'42.40"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "1.920.44"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").Value = "'6.292"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "'20.63"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "'7.310"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "'91.39"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "'18.02"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "'6.096"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "28.722.38"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "'2.279"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.139.39"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'162.76"
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.98"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.467"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'127.49"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1070"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.052"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.939"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.675"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02471"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.416"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06618"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2230"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6561"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.260"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.193"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "'5.008"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.13"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6173"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.31"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.753"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.302"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.014"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.237"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'122.62"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06954"
$ws.Range("E51").Value = "  +2.13%  "
